$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18: new match result added (04/08/2025 vs U. De Chile) ---
# Column A holds a date-looking string ("dd/mm/yyyy") that must stay plain
# text (inlineStr in the OOXML) rather than being auto-converted by Excel
# into a date serial number. Temporarily force the cell to Text format,
# assign the literal string, then restore the cell's style to Normal so
# no stray numeric/date formatting is left behind on the cell.
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "04/08/2025"
$ws.Range("A18").Style = "Normal"

$ws.Range("B18").Value = "U. De Chile"
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = "Cobresal"
$ws.Range("F18").Value = "W"
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0.38
$ws.Range("L18").Value = 1.68
$ws.Range("M18").Value = 7
$ws.Range("N18").Value = 22
$ws.Range("O18").Value = 4
$ws.Range("P18").Value = 3
